$wb = $excel.ActiveWorkbook
$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# Status column: "Ready for handoff" -> "Handed back: in sync with en-US"
$wsOverview.Range("E2").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F2").Value = "Handed back: in sync with en-US"
$wsZhCn.Range("C2").Value = "Handed back: in sync with en-US"
$wsDeDe.Range("C2").Value = "Handed back: in sync with en-US"

# zh-cn: Latest Handback DateTime updated, Error Detail cleared
$wsZhCn.Range("K2").Value = "2016-08-27 20:48:37"
$wsZhCn.Range("P2").Value = ""

# de-de: Latest Handback DateTime updated, Error Detail cleared
$wsDeDe.Range("K2").Value = "2016-08-27 20:48:44"
$wsDeDe.Range("P2").Value = ""

# Column width changes
$wsOverview.Range("E:F").ColumnWidth = 29.9777047293527
$wsZhCn.Range("C:C").ColumnWidth = 29.9777047293527
$wsZhCn.Range("P:P").ColumnWidth = 13.7470528738839
$wsDeDe.Range("C:C").ColumnWidth = 29.9777047293527
$wsDeDe.Range("P:P").ColumnWidth = 13.7470528738839
